$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Devices")
$ws.Activate()

# Update B3 with the new test case / method name (wraps to 2 lines, growing row 3's height)
$ws.Range("B3").Value = "VerifyAddUnitDetails"
$ws.Range("B3").Select()

# Rows 10-13: columns M and N change from text "NA" to boolean FALSE
$rows = 10,11,12,13
foreach ($r in $rows) {
    $ws.Cells.Item($r, 13).Value = $false
    $ws.Cells.Item($r, 14).Value = $false
}
